# Updates cryptocurrency price/volume data in the worksheet to match the
# latest scrape (GitHub Actions cryptos-list update).
#
# All target cells store values as text (inlineStr) even though many look
# like numbers (e.g. "1.00", "0.999"). Assigning such strings directly to
# Range.Value causes Excel's automatic type detection to coerce them into
# numeric values (losing formatting, e.g. "0.600" -> 0.6) and, for values
# containing special Unicode digits (e.g. the subscript-6 in "0.0₆0318"),
# can even corrupt the text entirely. To avoid this, every value is written
# through a temporary formula that evaluates to the exact literal string,
# then flattened back to a plain value via Copy / PasteSpecial (values only).
# This guarantees the text lands in the cell byte-for-byte, with no leftover
# formula and no cell-style changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Worksheet, $CellRef, $Text) {
    $range = $Worksheet.Range($CellRef)
    $escaped = $Text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $excel.CutCopyMode = $false
}

$updates = [ordered]@{
    'D2' = '68.614.52'
    'E2' = '  +0.68%  '
    'D3' = '2.700.29'
    'E3' = '  +2.18%  '
    'E4' = '  +0.00%  '
    'D5' = '599.08'
    'E5' = '  +0.39%  '
    'D6' = '160.55'
    'E6' = '  +2.90%  '
    'E7' = '  +0.02%  '
    'D8' = '0.544'
    'E8' = '  +0.51%  '
    'D9' = '2.700.01'
    'E9' = '  +2.20%  '
    'E10' = '  -0.38%  '
    'D11' = '0.156'
    'E11' = '  -0.41%  '
    'E12' = '  +1.18%  '
    'D13' = '0.361'
    'E13' = '  +2.66%  '
    'D14' = '28.28'
    'E14' = '  +1.05%  '
    'D15' = '3.190.94'
    'E15' = '  +2.24%  '
    'E16' = '  -0.47%  '
    'D17' = '68.526.73'
    'E17' = '  +0.63%  '
    'D18' = '2.712.73'
    'E18' = '  +2.91%  '
    'E19' = '  +4.35%  '
    'D20' = '367.51'
    'E20' = '  +1.25%  '
    'D21' = '7.65'
    'E21' = '  +3.89%  '
    'D22' = '4.54'
    'D23' = '4.90'
    'E23' = '  +2.37%  '
    'D24' = '2.13'
    'E24' = '  +3.32%  '
    'D25' = '74.50'
    'E25' = '  -0.51%  '
    'D26' = '1.00'
    'E26' = '  +0.03%  '
    'E27' = '  +3.06%  '
    'D28' = '2.825.64'
    'E28' = '  +1.87%  '
    'E29' = '  +0.91%  '
    'D30' = '0.999'
    'E30' = '  -14.42%  '
    'D31' = '578.12'
    'E31' = '  +4.50%  '
    'D32' = '8.26'
    'E32' = '  +2.94%  '
    'E33' = '  +3.61%  '
    'D34' = '1.95'
    'E34' = '  +5.90%  '
    'D35' = '0.133'
    'E35' = '  +3.85%  '
    'E36' = '  +6.58%  '
    'D38' = '161.63'
    'E38' = '  +0.21%  '
    'D39' = '19.91'
    'E39' = '  +1.62%  '
    'D40' = '0.380'
    'E40' = '  +2.21%  '
    'D41' = '1.92'
    'E41' = '  +2.75%  '
    'D42' = '5.41'
    'E42' = '  +1.76%  '
    'D43' = '2.70'
    'E43' = '  +3.71%  '
    'E44' = '  +0.40%  '
    'B45' = 'BabyDogeCoin'
    'C45' = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
    'D45' = '0.0₆0318'
    'E45' = '  -5.28%  '
    'B46' = 'USDe'
    'C46' = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
    'D46' = '1.00'
    'E46' = '  +0.05%  '
    'D47' = '157.78'
    'E47' = '  -0.81%  '
    'E48' = '  +6.80%  '
    'E49' = '  +5.02%  '
    'D50' = '0.600'
    'E50' = '  +6.89%  '
    'D51' = '22.05'
    'E51' = '  -0.27%  '
}

foreach ($cellRef in $updates.Keys) {
    Set-TextValue $ws $cellRef $updates[$cellRef]
}

